$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; New="53÷7="},
    @{Row=1;  Col=2; New="54÷2="},
    @{Row=1;  Col=3; New="98÷3="},
    @{Row=1;  Col=4; New="69÷7="},
    @{Row=1;  Col=5; New="89÷2="},
    @{Row=5;  Col=1; New="90÷7="},
    @{Row=5;  Col=2; New="89÷6="},
    @{Row=5;  Col=3; New="87÷5="},
    @{Row=5;  Col=4; New="90÷2="},
    @{Row=5;  Col=5; New="39÷5="},
    @{Row=9;  Col=1; New="61÷5="},
    @{Row=9;  Col=2; New="52÷4="},
    @{Row=9;  Col=3; New="66÷8="},
    @{Row=9;  Col=4; New="52÷2="},
    @{Row=9;  Col=5; New="58÷7="},
    @{Row=13; Col=1; New="79÷7="},
    @{Row=13; Col=2; New="77÷2="},
    @{Row=13; Col=3; New="68÷8="},
    @{Row=13; Col=4; New="74÷9="},
    @{Row=13; Col=5; New="70÷9="},
    @{Row=17; Col=1; New="83÷4="},
    @{Row=17; Col=2; New="19÷3="},
    @{Row=17; Col=3; New="92÷6="},
    @{Row=17; Col=4; New="35÷6="},
    @{Row=17; Col=5; New="66÷4="}
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $rng = $cell.Range
    $rng.Text = $item.New
}
